# Conclusion slide (slide 15): grow the body placeholder and add a new
# "Delete Account" bullet under the "Future releases include-" list.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(2)

# 1) Grow the content placeholder's height (8281350 x 4314324 -> 8281350 x 5586162 EMU).
#    PowerPoint's Shape.Height/Top/Width/Left are expressed in points (1 pt = 12700 EMU).
$shp.Height = 5586162 / 12700

# 2) Append a new level-1 bullet paragraph ("Delete Account") after the last
#    paragraph ("Using Jtoken for validation"), inheriting that list's bullet
#    (Wingdings "Ø") and matching the run formatting used throughout this box.
$tr = $shp.TextFrame.TextRange
$startLen = $tr.Length

$tr.InsertAfter([char]13 + "Delete Account") | Out-Null

$newRun = $tr.Characters($startLen + 2, 14)
$newRun.Font.Size = 18
$newRun.Font.Bold = 0
$newRun.Font.Italic = 0
$newRun.Font.Underline = 0
$newRun.Font.Strikethrough = 0
$newRun.Font.BaselineOffset = 0
